$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a brand-new item row ("FLACORT 6MG 20 TABS.") above what is
# currently row 21 ("FLUB 20MG/ML SUSPENSION 30 ML"). Everything currently
# at/after row 21 (items, the totals row, and the footer row) shifts down by
# one row.
# ---------------------------------------------------------------------------
$ws.Rows("21:21").Insert()

# Copy the formatting (styles/borders/number-formats) of a normal item row
# (row 22, which after the shift carries the old row-21 formatting) onto the
# freshly inserted - currently blank - row 21.
$ws.Range("A22:Q22").Copy()
$ws.Range("A21:Q21").PasteSpecial(-4122)

# Recreate the merged cells that every item row has (A:B, C:G, H:K, L:M, N:O)
# for the new row 21 - Insert() does not clone merges automatically.
$ws.Range("A21:B21").Merge()
$ws.Range("C21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("N21:O21").Merge()

# Fill in the new row's data (م / الاسم / الرصيد الحالي / حد الطلب / السعر / سعر البيع / عدد التعاملات).
$ws.Range("A21").Value2 = 15
$ws.Range("C21").Value2 = "FLACORT 6MG 20 TABS."
$ws.Range("H21").Value2 = "1:0"
$ws.Range("L21").Value2 = "1"
$ws.Range("N21").Value2 = "54.00"
$ws.Range("P21").Value2 = "27.0000"
$ws.Range("Q21").Value2 = "0:1"

# The serial-number column (م) is independent of the row shift: it must stay
# the simple sequence 16,17,18,...,35 for the rows that used to be 21-40 and
# are now 22-41 (Insert() dragged the old values - 15,16,...,34 - down with
# the rows, so they are all off by one and need correcting).
for ($r = 22; $r -le 41; $r++) {
  $ws.Range("A$r").Value2 = $r - 6
}

# Re-apply the exact row heights the finished report uses. Rows 21-40 keep
# the height that already belonged to that row number; row 41 is the new
# last item row; rows 42/43 are the shifted totals/footer rows.
$ws.Rows("21:21").RowHeight = 25.5
$ws.Rows("22:22").RowHeight = 25.5
$ws.Rows("23:23").RowHeight = 24.75
$ws.Rows("24:24").RowHeight = 25.5
$ws.Rows("25:25").RowHeight = 24.75
$ws.Rows("26:26").RowHeight = 25.5
$ws.Rows("27:27").RowHeight = 25.5
$ws.Rows("28:28").RowHeight = 24.75
$ws.Rows("29:29").RowHeight = 25.5
$ws.Rows("30:30").RowHeight = 24.75
$ws.Rows("31:31").RowHeight = 25.5
$ws.Rows("32:32").RowHeight = 25.5
$ws.Rows("33:33").RowHeight = 24.75
$ws.Rows("34:34").RowHeight = 25.5
$ws.Rows("35:35").RowHeight = 24.75
$ws.Rows("36:36").RowHeight = 25.5
$ws.Rows("37:37").RowHeight = 25.5
$ws.Rows("38:38").RowHeight = 24.75
$ws.Rows("39:39").RowHeight = 25.5
$ws.Rows("40:40").RowHeight = 24.75
$ws.Rows("41:41").RowHeight = 25.5
$ws.Rows("42:42").RowHeight = 25.5
$ws.Rows("43:43").RowHeight = 16.5

# ---------------------------------------------------------------------------
# Update the "Total" row (now row 42): add the new item's sell price
# (27.0000) to the previous total (2041.54).
# ---------------------------------------------------------------------------
$ws.Range("P42").Value2 = 2068.54

# ---------------------------------------------------------------------------
# Update the generated-on timestamp in the footer (now row 43) to the new
# save time.
# ---------------------------------------------------------------------------
$ws.Range("A43").Value2 = "Tuesday, 19 August, 2025 12:49 PM"

Write-Host "Applied FLACORT row insert + timestamp update"
